$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: fill in F6 (previously empty) ---
$ws.Range("F6").Value = "1.Clic en botón ""Editar""`n2.Diligenciar campo ""Observaciones"" en modal Editar`n3.Clic en botón ""Editar / Guardar"" en modal`n4.Cerrar modal de edición`n"

# --- Row 7: new test case row, mirroring row 6's formatting ---
$ws.Rows("6").Copy()
$ws.Rows("7").PasteSpecial()

$ws.Range("A7").Value = "CP_INFTECRED_006"
$ws.Range("B7").Value = "refrescar"
$ws.Range("C7").Value = "Positivo"
$ws.Range("D7").Value = "eCenter"
$ws.Range("E7").Value = "Haber ingresado a la vista"
$ws.Range("F7").Value = "1.Clic en botón ""Refrescar"""
$ws.Range("G7").Value = "N/A"
$ws.Range("H7").Value = "El sistema dede refrescar la vista correctamente"
$ws.Range("I7").Value = "la vista se refresca correctamente"
$ws.Range("J7").Value = "OK"
$ws.Range("K7").Value = "SI"
$ws.Range("L7").Value = "N/A"
